# Insert a new weekly data row for "Cebollín" (Terminal Hortofrutícola Agro Chillán)
# before the current row 179. This pushes the existing rows 179-214 down to 180-215
# (dimension grows from A1:R214 to A1:R215), and fills the newly inserted row 179
# with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 179, shifting rows 179:214 down to 180:215.
$ws.Rows.Item(179).Insert()

$row = 179

$ws.Cells.Item($row, 1).Value2  = 7
$ws.Cells.Item($row, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value2  = "Ñuble"
$ws.Cells.Item($row, 4).Value2  = 45204
$ws.Cells.Item($row, 5).Value2  = 16
$ws.Cells.Item($row, 6).Value2  = 100112037
$ws.Cells.Item($row, 7).Value2  = "Cebollín"
$ws.Cells.Item($row, 8).Value2  = "Sin especificar"
$ws.Cells.Item($row, 9).Value2  = "Primera"
$ws.Cells.Item($row, 10).Value2 = 250
$ws.Cells.Item($row, 11).Value2 = 6000
$ws.Cells.Item($row, 12).Value2 = 6000
$ws.Cells.Item($row, 13).Value2 = 6000
$ws.Cells.Item($row, 14).Value2 = "`$/paquete 36 unidades"
$ws.Cells.Item($row, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item($row, 16).Value2 = 167
$ws.Cells.Item($row, 17).Value2 = 36
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"

# Preserve the date number format (yyyy-mm-dd hh:mm:ss) used by the other rows
# in column D for the newly inserted row.
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row + 1, 4).NumberFormat
